$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 0

$ws.Range("D9").Select()
